$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.292.61'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '1.647.12'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").Value = "'217.40"

$ws.Range("E6").Value = '  +0.42%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").Value = "'19.96"
$ws.Range("E10").Value = '  +0.90%  '

$ws.Range("E11").Value = '  +0.07%  '

$ws.Range("E12").Value = '  +0.54%  '

$ws.Range("D13").Value = '1.874.92'
$ws.Range("E13").Value = '  +0.40%  '

$ws.Range("D14").Value = '1.664.79'
$ws.Range("E14").Value = '  +1.45%  '

$ws.Range("E15").Value = '  -2.39%  '

$ws.Range("E16").Value = '  -0.38%  '

$ws.Range("D17").Value = "'63.54"
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").Value = '26.267.62'
$ws.Range("E18").Value = '  +1.49%  '

$ws.Range("E19").Value = '  -0.12%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = "'195.98"
$ws.Range("E20").Value = '  +1.36%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'4.45"
$ws.Range("E21").Value = '  -0.84%  '

$ws.Range("D22").Value = "'10.07"
$ws.Range("E22").Value = '  +0.71%  '

$ws.Range("E23").Value = '  -0.48%  '

$ws.Range("D24").Value = "'143.27"
$ws.Range("E24").Value = '  +0.38%  '

$ws.Range("E25").Value = '  -2.56%  '

$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("E27").Value = '  +2.16%  '

$ws.Range("E28").Value = '  -0.11%  '

$ws.Range("D29").Value = "'15.63"
$ws.Range("E29").Value = '  +0.28%  '

$ws.Range("E30").Value = '  +1.27%  '

$ws.Range("D31").Value = "'0.0507"
$ws.Range("E31").Value = '  +2.00%  '

$ws.Range("E32").Value = '  -0.24%  '

$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = '  -0.17%  '

$ws.Range("E34").Value = '  +1.57%  '

$ws.Range("E35").Value = '  +1.03%  '

$ws.Range("E36").Value = '  +0.28%  '

$ws.Range("D37").Value = '1.138.09'
$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("D38").Value = "'0.555"
$ws.Range("E38").Value = '  +1.45%  '

$ws.Range("E39").Value = '  -1.36%  '

$ws.Range("E40").Value = '  +0.32%  '

$ws.Range("E41").Value = '  -0.15%  '

$ws.Range("D42").Value = "'5.70"
$ws.Range("E42").Value = '  +2.47%  '

$ws.Range("D43").Value = "'100.34"
$ws.Range("E43").Value = '  -0.33%  '

$ws.Range("D44").Value = "'0.799"
$ws.Range("E44").Value = '  -1.38%  '

$ws.Range("D45").Value = '1.783.61'
$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("D46").Value = "'56.33"
$ws.Range("E46").Value = '  +1.44%  '

$ws.Range("E47").Value = '  +2.90%  '

$ws.Range("D48").Value = "'0.0514"
$ws.Range("E48").Value = '  +1.99%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'7.72"
$ws.Range("E49").Value = '  +2.47%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = "'0.418"
$ws.Range("E50").Value = '  -0.12%  '

$ws.Range("D51").Value = "'0.0974"
$ws.Range("E51").Value = '  +1.69%  '
